# إضافة حدث جديد في Card12 by admin at 2025-12-08 11:54:20
#
# Card12's service-log sheet gains a new row (row 24) recording a fresh
# maintenance event, and the previously-blank tracking columns on the prior
# row (23) get backfilled with the sheet's usual "nan" placeholder text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Card12")

# Row 23: the measurement columns (B..K) and the Event column (N) were left
# completely empty; fill them with the literal placeholder "nan" used
# throughout the rest of the sheet for not-applicable values.
$ws.Cells.Item(23, 2).Value = "nan"
$ws.Cells.Item(23, 3).Value = "nan"
$ws.Cells.Item(23, 4).Value = "nan"
$ws.Cells.Item(23, 5).Value = "nan"
$ws.Cells.Item(23, 6).Value = "nan"
$ws.Cells.Item(23, 7).Value = "nan"
$ws.Cells.Item(23, 8).Value = "nan"
$ws.Cells.Item(23, 9).Value = "nan"
$ws.Cells.Item(23, 10).Value = "nan"
$ws.Cells.Item(23, 11).Value = "nan"
$ws.Cells.Item(23, 14).Value = "nan"

# Row 24: brand-new service event for card 12.
#   A = card number (kept as text, matching every other cell in column A)
#   L = Date, M = Serviced by, O = Correction note
$ws.Cells.Item(24, 1).NumberFormat = "@"
$ws.Cells.Item(24, 1).Value = "12"
$ws.Cells.Item(24, 12).Value = "30\1\2025"
$ws.Cells.Item(24, 13).Value = "الخبير"
$ws.Cells.Item(24, 15).Value = "تم سن الفلاتس وتغيير الجرائد الخلفيه (1_5_8)"

Write-Output "Card12: backfilled row 23 and appended row 24"
